# Update cryptos list values (price / volume %) per the latest scrape.
# Columns: A=rank(idx), B=Coin, C=Link, D=Price, E=Volume(1h)
# D/E columns are stored as text in the workbook (e.g. "317.29", "  -2.62%  "),
# so price values that look numeric must be written with a Text number format
# to stop Excel from auto-converting them to real numbers; the format is then
# restored to the default "Normal" style so no formatting diff is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Row 2 - Bitcoin
Set-TextCell 2 4 "43.262.08"
Set-Cell     2 5 "  -2.23%  "

# Row 3 - Ethereum
Set-TextCell 3 4 "2.228.57"
Set-Cell     3 5 "  -4.25%  "

# Row 4 - TetherUSD
Set-Cell 4 5 "  -0.14%  "

# Row 5 - BNB
Set-TextCell 5 4 "317.29"
Set-Cell     5 5 "  -2.62%  "

# Row 6 - Solana
Set-TextCell 6 4 "98.61"
Set-Cell     6 5 "  -5.46%  "

# Row 7 - XRP
Set-TextCell 7 4 "0.581"
Set-Cell     7 5 "  -6.99%  "

# Row 8 - USDC
Set-Cell 8 5 "  -0.16%  "

# Row 9 - Cardano
Set-TextCell 9 4 "0.560"
Set-Cell     9 5 "  -6.92%  "

# Row 10 - Avalanche
Set-TextCell 10 4 "37.07"
Set-Cell     10 5 "  -6.82%  "

# Row 11 - OKB
Set-TextCell 11 4 "53.99"
Set-Cell     11 5 "  -2.35%  "

# Row 12 - Dogecoin
Set-Cell 12 5 "  -8.27%  "

# Row 13 - Polkadot
Set-Cell 13 5 "  -6.42%  "

# Row 14 - TRON
Set-Cell 14 5 "  -1.15%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextCell 15 4 "2.567.13"
Set-Cell     15 5 "  -4.37%  "

# Row 16 - Polygon
Set-TextCell 16 4 "0.861"
Set-Cell     16 5 "  -10.08%  "

# Row 17 - Chainlink
Set-TextCell 17 4 "14.35"
Set-Cell     17 5 "  -4.80%  "

# Row 18 - WrappedEther
Set-TextCell 18 4 "2.231.56"
Set-Cell     18 5 "  -5.01%  "

# Row 19 - WrappedBTC
Set-TextCell 19 4 "43.184.66"
Set-Cell     19 5 "  -2.55%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextCell 20 4 "14.01"
Set-Cell     20 5 "  -7.08%  "

# Row 21 - ShibaInu
Set-TextCell 21 4 "0.0₃0966"
Set-Cell     21 5 "  -7.46%  "

# Row 22 - Uniswap
Set-Cell 22 5 "  -8.36%  "

# Row 23 - PancakeSwap
Set-Cell 23 5 "  -9.45%  "

# Row 24 - Litecoin
Set-TextCell 24 4 "65.14"
Set-Cell     24 5 "  -9.98%  "

# Row 25 - BitcoinCash
Set-TextCell 25 4 "236.93"
Set-Cell     25 5 "  -6.82%  "

# Row 26 - ImmutableX
Set-TextCell 26 4 "2.19"
Set-Cell     26 5 "  -2.01%  "

# Row 27 - Dai
Set-Cell 27 5 "  +0.32%  "

# Row 28 - LEO
Set-Cell 28 5 "  +2.18%  "

# Row 29 - Cosmos
Set-TextCell 29 4 "10.05"
Set-Cell     29 5 "  -8.90%  "

# Row 30 - Toncoin
Set-TextCell 30 4 "2.22"
Set-Cell     30 5 "  -2.44%  "

# Row 31 - Filecoin
Set-TextCell 31 4 "6.39"
Set-Cell     31 5 "  -11.47%  "

# Row 32 - InjectiveProtocol
Set-TextCell 32 4 "36.58"
Set-Cell     32 5 "  +0.84%  "

# Row 33 - EthereumClassic
Set-TextCell 33 4 "20.26"
Set-Cell     33 5 "  -6.62%  "

# Row 34 - now Hedera (was Monero)
Set-Cell     34 2 "Hedera"
Set-Cell     34 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 34 4 "0.0866"
Set-Cell     34 5 "  -7.49%  "

# Row 35 - now Monero (was Hedera)
Set-Cell     35 2 "Monero"
Set-Cell     35 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 35 4 "157.68"
Set-Cell     35 5 "  -4.10%  "

# Row 36 - now LidoDAOToken (was WEMIXToken)
Set-Cell     36 2 "LidoDAOToken"
Set-Cell     36 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 36 4 "3.31"
Set-Cell     36 5 "  +4.05%  "

# Row 37 - now WEMIXToken (was LidoDAOToken)
Set-Cell     37 2 "WEMIXToken"
Set-Cell     37 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell 37 4 "2.67"
Set-Cell     37 5 "  -4.10%  "

# Row 38 - Stellar
Set-TextCell 38 4 "0.121"
Set-Cell     38 5 "  -6.81%  "

# Row 39 - RenderToken
Set-TextCell 39 4 "4.40"
Set-Cell     39 5 "  -4.24%  "

# Row 40 - ARBITRUM
Set-Cell 40 5 "  -2.41%  "

# Row 41 - Kaspa
Set-Cell 41 5 "  -7.82%  "

# Row 42 - NEARProtocol
Set-Cell 42 5 "  -3.68%  "

# Row 43 - VeChain
Set-Cell 43 5 "  -7.80%  "

# Row 44 - Celestia
Set-TextCell 44 4 "14.40"
Set-Cell     44 5 "  +14.29%  "

# Row 45 - FirstDigitalUSD
Set-Cell 45 5 "  -0.19%  "

# Row 46 - Maker
Set-TextCell 46 4 "1.767.73"
Set-Cell     46 5 "  -5.05%  "

# Row 47 - Algorand
Set-TextCell 47 4 "0.203"
Set-Cell     47 5 "  -8.92%  "

# Row 48 - BitcoinSV
Set-TextCell 48 4 "83.73"
Set-Cell     48 5 "  -10.48%  "

# Row 49 - FraxShare
Set-TextCell 49 4 "8.88"
Set-Cell     49 5 "  -2.38%  "

# Row 50 - THORChain
Set-TextCell 50 4 "5.28"
Set-Cell     50 5 "  -10.05%  "

# Row 51 - ordi
Set-TextCell 51 4 "73.49"
Set-Cell     51 5 "  -9.68%  "
